$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 56.756364
$ws.Range("H2").Value = 170.269092
$ws.Range("I2").Value = 0.157357217290148
$ws.Range("J2").Value = 0.157357217290148
$ws.Range("M2").Value = 6.493001333333333
$ws.Range("N2").Value = 19.479004
$ws.Range("O2").Value = 0.01042978512569556
$ws.Range("P2").Value = 0.01042978512569556
$ws.Range("Q2").Value = 368.5191471271519
$ws.Range("R2").Value = 3316.672324144368
$ws.Range("S2").Value = 0.00164120196431363
$ws.Range("T2").Value = 0.00164120196431363
$ws.Range("G3").Value = 56.756364
$ws.Range("H3").Value = 170.269092
$ws.Range("I3").Value = 0.157357217290148
$ws.Range("J3").Value = 0.157357217290148
$ws.Range("O3").Value = 0.03234404904340005
$ws.Range("P3").Value = 0.03234404904340005
$ws.Range("Q3").Value = 1142.823291608096
$ws.Range("R3").Value = 10285.40962447286
$ws.Range("S3").Value = 0.005089569553365503
$ws.Range("T3").Value = 0.005089569553365504
$ws.Range("G4").Value = 56.756364
$ws.Range("H4").Value = 170.269092
$ws.Range("I4").Value = 0.157357217290148
$ws.Range("J4").Value = 0.157357217290148
$ws.Range("M4").Value = 300.2397663333333
$ws.Range("N4").Value = 900.719299
$ws.Range("O4").Value = 0.4822787010638293
$ws.Range("P4").Value = 0.4822787010638293
$ws.Range("Q4").Value = 17040.51746528961
$ws.Range("R4").Value = 153364.6571876065
$ws.Range("S4").Value = 0.07589003435771129
$ws.Range("T4").Value = 0.0758900343577113
$ws.Range("G5").Value = 56.756364
$ws.Range("H5").Value = 170.269092
$ws.Range("I5").Value = 0.157357217290148
$ws.Range("J5").Value = 0.157357217290148
$ws.Range("M5").Value = 2.838981666666667
$ws.Range("N5").Value = 8.516945
$ws.Range("O5").Value = 0.004560289955141813
$ws.Range("P5").Value = 0.004560289955141813
$ws.Range("Q5").Value = 161.13027686266
$ws.Range("R5").Value = 1450.17249176394
$ws.Range("S5").Value = 0.0007175945373773294
$ws.Range("T5").Value = 0.0007175945373773295
$ws.Range("G6").Value = 56.756364
$ws.Range("H6").Value = 170.269092
$ws.Range("I6").Value = 0.157357217290148
$ws.Range("J6").Value = 0.157357217290148
$ws.Range("M6").Value = 292.8367666666667
$ws.Range("N6").Value = 878.5103
$ws.Range("O6").Value = 0.4703871748119333
$ws.Range("P6").Value = 0.4703871748119333
$ws.Range("Q6").Value = 16620.3501215164
$ws.Range("R6").Value = 149583.1510936476
$ws.Range("S6").Value = 0.0740188168773802
$ws.Range("T6").Value = 0.07401881687738021
$ws.Range("I7").Value = 0.1972445753159741
$ws.Range("J7").Value = 0.1972445753159741
$ws.Range("M7").Value = 6.493001333333333
$ws.Range("N7").Value = 19.479004
$ws.Range("O7").Value = 0.01042978512569556
$ws.Range("P7").Value = 0.01042978512569556
$ws.Range("Q7").Value = 461.9324357831729
$ws.Range("R7").Value = 4157.391922048556
$ws.Range("S7").Value = 0.002057218537754684
$ws.Range("T7").Value = 0.002057218537754685
$ws.Range("I8").Value = 0.1972445753159741
$ws.Range("J8").Value = 0.1972445753159741
$ws.Range("O8").Value = 0.03234404904340005
$ws.Range("P8").Value = 0.03234404904340005
$ws.Range("S8").Value = 0.006379688217564479
$ws.Range("T8").Value = 0.00637968821756448
$ws.Range("I9").Value = 0.1972445753159741
$ws.Range("J9").Value = 0.1972445753159741
$ws.Range("M9").Value = 300.2397663333333
$ws.Range("N9").Value = 900.719299
$ws.Range("O9").Value = 0.4822787010638293
$ws.Range("P9").Value = 0.4822787010638293
$ws.Range("Q9").Value = 21359.99662734203
$ws.Range("R9").Value = 192239.9696460783
$ws.Range("S9").Value = 0.0951268575752746
$ws.Range("T9").Value = 0.09512685757527461
$ws.Range("I10").Value = 0.1972445753159741
$ws.Range("J10").Value = 0.1972445753159741
$ws.Range("M10").Value = 2.838981666666667
$ws.Range("N10").Value = 8.516945
$ws.Range("O10").Value = 0.004560289955141813
$ws.Range("P10").Value = 0.004560289955141813
$ws.Range("Q10").Value = 201.9740408329561
$ws.Range("R10").Value = 1817.766367496605
$ws.Range("S10").Value = 0.0008994924555196493
$ws.Range("T10").Value = 0.0008994924555196495
$ws.Range("I11").Value = 0.1972445753159741
$ws.Range("J11").Value = 0.1972445753159741
$ws.Range("M11").Value = 292.8367666666667
$ws.Range("N11").Value = 878.5103
$ws.Range("O11").Value = 0.4703871748119333
$ws.Range("P11").Value = 0.4703871748119333
$ws.Range("Q11").Value = 20833.32406213408
$ws.Range("R11").Value = 187499.9165592067
$ws.Range("S11").Value = 0.09278131852986063
$ws.Range("T11").Value = 0.09278131852986064
$ws.Range("G12").Value = 124.7878343333334
$ws.Range("H12").Value = 374.363503
$ws.Range("I12").Value = 0.3459747062436438
$ws.Range("J12").Value = 0.3459747062436438
$ws.Range("M12").Value = 6.493001333333333
$ws.Range("N12").Value = 19.479004
$ws.Range("O12").Value = 0.01042978512569556
$ws.Range("P12").Value = 0.01042978512569556
$ws.Range("Q12").Value = 810.2475747101125
$ws.Range("R12").Value = 7292.228172391013
$ws.Range("S12").Value = 0.003608441845046848
$ws.Range("T12").Value = 0.003608441845046849
$ws.Range("G13").Value = 124.7878343333334
$ws.Range("H13").Value = 374.363503
$ws.Range("I13").Value = 0.3459747062436438
$ws.Range("J13").Value = 0.3459747062436438
$ws.Range("O13").Value = 0.03234404904340005
$ws.Range("P13").Value = 0.03234404904340005
$ws.Range("Q13").Value = 2512.677584234709
$ws.Range("R13").Value = 22614.09825811238
$ws.Range("S13").Value = 0.01119022286652034
$ws.Range("T13").Value = 0.01119022286652034
$ws.Range("G14").Value = 124.7878343333334
$ws.Range("H14").Value = 374.363503
$ws.Range("I14").Value = 0.3459747062436438
$ws.Range("J14").Value = 0.3459747062436438
$ws.Range("M14").Value = 300.2397663333333
$ws.Range("N14").Value = 900.719299
$ws.Range("O14").Value = 0.4822787010638293
$ws.Range("P14").Value = 0.4822787010638293
$ws.Range("Q14").Value = 37466.27022148271
$ws.Range("R14").Value = 337196.4319933444
$ws.Range("S14").Value = 0.1668562319281245
$ws.Range("T14").Value = 0.1668562319281245
$ws.Range("G15").Value = 124.7878343333334
$ws.Range("H15").Value = 374.363503
$ws.Range("I15").Value = 0.3459747062436438
$ws.Range("J15").Value = 0.3459747062436438
$ws.Range("M15").Value = 2.838981666666667
$ws.Range("N15").Value = 8.516945
$ws.Range("O15").Value = 0.004560289955141813
$ws.Range("P15").Value = 0.004560289955141813
$ws.Range("Q15").Value = 354.2703738953706
$ws.Range("R15").Value = 3188.433365058335
$ws.Range("S15").Value = 0.001577744977616029
$ws.Range("T15").Value = 0.001577744977616029
$ws.Range("G16").Value = 124.7878343333334
$ws.Range("H16").Value = 374.363503
$ws.Range("I16").Value = 0.3459747062436438
$ws.Range("J16").Value = 0.3459747062436438
$ws.Range("M16").Value = 292.8367666666667
$ws.Range("N16").Value = 878.5103
$ws.Range("O16").Value = 0.4703871748119333
$ws.Range("P16").Value = 0.4703871748119333
$ws.Range("Q16").Value = 36542.46592550899
$ws.Range("R16").Value = 328882.1933295809
$ws.Range("S16").Value = 0.1627420646263362
$ws.Range("T16").Value = 0.1627420646263362
$ws.Range("G17").Value = 21.495283
$ws.Range("H17").Value = 64.485849
$ws.Range("I17").Value = 0.05959574714377799
$ws.Range("J17").Value = 0.05959574714377799
$ws.Range("M17").Value = 6.493001333333333
$ws.Range("N17").Value = 19.479004
$ws.Range("O17").Value = 0.01042978512569556
$ws.Range("P17").Value = 0.01042978512569556
$ws.Range("Q17").Value = 139.5689011793773
$ws.Range("R17").Value = 1256.120110614396
$ws.Range("S17").Value = 0.0006215708371148895
$ws.Range("T17").Value = 0.0006215708371148896
$ws.Range("G18").Value = 21.495283
$ws.Range("H18").Value = 64.485849
$ws.Range("I18").Value = 0.05959574714377799
$ws.Range("J18").Value = 0.05959574714377799
$ws.Range("O18").Value = 0.03234404904340005
$ws.Range("P18").Value = 0.03234404904340005
$ws.Range("Q18").Value = 432.8203630540453
$ws.Range("R18").Value = 3895.383267486408
$ws.Range("S18").Value = 0.001927567768396424
$ws.Range("T18").Value = 0.001927567768396424
$ws.Range("G19").Value = 21.495283
$ws.Range("H19").Value = 64.485849
$ws.Range("I19").Value = 0.05959574714377799
$ws.Range("J19").Value = 0.05959574714377799
$ws.Range("M19").Value = 300.2397663333333
$ws.Range("N19").Value = 900.719299
$ws.Range("O19").Value = 0.4822787010638293
$ws.Range("P19").Value = 0.4822787010638293
$ws.Range("Q19").Value = 6453.738745188872
$ws.Range("R19").Value = 58083.64870669985
$ws.Range("S19").Value = 0.02874175952142966
$ws.Range("T19").Value = 0.02874175952142966
$ws.Range("G20").Value = 21.495283
$ws.Range("H20").Value = 64.485849
$ws.Range("I20").Value = 0.05959574714377799
$ws.Range("J20").Value = 0.05959574714377799
$ws.Range("M20").Value = 2.838981666666667
$ws.Range("N20").Value = 8.516945
$ws.Range("O20").Value = 0.004560289955141813
$ws.Range("P20").Value = 0.004560289955141813
$ws.Range("Q20").Value = 61.02471435681167
$ws.Range("R20").Value = 549.222429211305
$ws.Range("S20").Value = 0.0002717738870689422
$ws.Range("T20").Value = 0.0002717738870689422
$ws.Range("G21").Value = 21.495283
$ws.Range("H21").Value = 64.485849
$ws.Range("I21").Value = 0.05959574714377799
$ws.Range("J21").Value = 0.05959574714377799
$ws.Range("M21").Value = 292.8367666666667
$ws.Range("N21").Value = 878.5103
$ws.Range("O21").Value = 0.4703871748119333
$ws.Range("P21").Value = 0.4703871748119333
$ws.Range("Q21").Value = 6294.609172304967
$ws.Range("R21").Value = 56651.4825507447
$ws.Range("S21").Value = 0.02803307512976807
$ws.Range("T21").Value = 0.02803307512976807
$ws.Range("G22").Value = 86.50223699999999
$ws.Range("H22").Value = 259.506711
$ws.Range("I22").Value = 0.239827754006456
$ws.Range("J22").Value = 0.2398277540064561
$ws.Range("M22").Value = 6.493001333333333
$ws.Range("N22").Value = 19.479004
$ws.Range("O22").Value = 0.01042978512569556
$ws.Range("P22").Value = 0.01042978512569556
$ws.Range("Q22").Value = 561.659140177316
$ws.Range("R22").Value = 5054.932261595844
$ws.Range("S22").Value = 0.002501351941465509
$ws.Range("T22").Value = 0.00250135194146551
$ws.Range("G23").Value = 86.50223699999999
$ws.Range("H23").Value = 259.506711
$ws.Range("I23").Value = 0.239827754006456
$ws.Range("J23").Value = 0.2398277540064561
$ws.Range("O23").Value = 0.03234404904340005
$ws.Range("P23").Value = 0.03234404904340005
$ws.Range("Q23").Value = 1741.774212664568
$ws.Range("R23").Value = 15675.96791398111
$ws.Range("S23").Value = 0.007757000637553296
$ws.Range("T23").Value = 0.007757000637553297
$ws.Range("G24").Value = 86.50223699999999
$ws.Range("H24").Value = 259.506711
$ws.Range("I24").Value = 0.239827754006456
$ws.Range("J24").Value = 0.2398277540064561
$ws.Range("M24").Value = 300.2397663333333
$ws.Range("N24").Value = 900.719299
$ws.Range("O24").Value = 0.4822787010638293
$ws.Range("P24").Value = 0.4822787010638293
$ws.Range("Q24").Value = 25971.41142419062
$ws.Range("R24").Value = 233742.7028177156
$ws.Range("S24").Value = 0.1156638176812892
$ws.Range("T24").Value = 0.1156638176812892
$ws.Range("G25").Value = 86.50223699999999
$ws.Range("H25").Value = 259.506711
$ws.Range("I25").Value = 0.239827754006456
$ws.Range("J25").Value = 0.2398277540064561
$ws.Range("M25").Value = 2.838981666666667
$ws.Range("N25").Value = 8.516945
$ws.Range("O25").Value = 0.004560289955141813
$ws.Range("P25").Value = 0.004560289955141813
$ws.Range("Q25").Value = 245.578264968655
$ws.Range("R25").Value = 2210.204384717895
$ws.Range("S25").Value = 0.001093684097559863
$ws.Range("T25").Value = 0.001093684097559863
$ws.Range("G26").Value = 86.50223699999999
$ws.Range("H26").Value = 259.506711
$ws.Range("I26").Value = 0.239827754006456
$ws.Range("J26").Value = 0.2398277540064561
$ws.Range("M26").Value = 292.8367666666667
$ws.Range("N26").Value = 878.5103
$ws.Range("O26").Value = 0.4703871748119333
$ws.Range("P26").Value = 0.4703871748119333
$ws.Range("Q26").Value = 25331.0353925137
$ws.Range("R26").Value = 227979.3185326233
$ws.Range("S26").Value = 0.1128118996485882
$ws.Range("T26").Value = 0.1128118996485882
